$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for season record columns, matching the style
# already used by the other header cells (bold, bordered, centered).
$ws.Range("AA1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record for every player row (2-46) with the
# team's win/loss/tie totals for the season.
$ws.Range("AD2:AD46").Value = 78
$ws.Range("AE2:AE46").Value = 83
$ws.Range("AF2:AF46").Value = 0

Write-Host "Season record columns added"
